$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "ABC-1234" to "Teste-1234" (also updates the
# _FilterDatabase defined name reference automatically)
$ws.Name = "Teste-1234"

# Change selection from A1:A1048576 to active cell B21
$ws.Range("B21").Select()
